# Spez preis 21.4.24 plus Ella 2. Vorstelung ohne Mindestabgaben
# Adds a new row (13) to the "Spezialpreisekiosk" table with the special
# price entry for 21.4.2024 ("Treppensitze").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / styles) from the last existing
# data row (12) down to the new row (13) for the date and price columns,
# so the new cells reuse the existing style indices instead of creating
# new ones.
$ws.Range("A12").Copy($ws.Range("A13")) | Out-Null
$ws.Range("D12").Copy($ws.Range("D13")) | Out-Null

# Fill in the new row's values.
$ws.Range("A13").Value = 45403
$ws.Range("B13").Value = "Spez 1"
$ws.Range("C13").Value = "Treppensitze"
$ws.Range("D13").Value = 10

# Grow the worksheet table so it includes the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E13")) | Out-Null

# Match the active selection left behind after entering the new row.
$ws.Range("C14").Select() | Out-Null
